$d = $word.ActiveDocument

# Locate the existing "Mandarin" run and collapse the found range to its end
$rng = $d.Content
[void]$rng.Find.Execute("Mandarin", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$rng.Collapse(0)

# Insert the new text as its own run right after "Mandarin"
$rng.InsertAfter(" (Simplified)")

# Touch the font property (set then restore) so the new text becomes a
# distinct run instead of being merged back into the preceding "Mandarin" run
$rng.Font.Bold = 1
$rng.Font.Bold = 0
